# Apply the authoring changes described in the commit:
#  - "updated ubc2 10 scripts": the Neo4j query stored in cell B2 of the
#    "startup" sheet (CasesTab query) had its trailing `Cohort` column
#    removed from the RETURN clause.
#  - the workbook was left with cell B2 selected (and scrolled so row 2
#    is visible) instead of D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldQuery = $ws.Range("B2").Value()

$cohortSuffix = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"

if ($oldQuery.EndsWith($cohortSuffix)) {
    $newQuery = $oldQuery.Substring(0, $oldQuery.Length - $cohortSuffix.Length)
    $ws.Range("B2").Value = $newQuery
}

# Reflect the saved selection/view state: active cell moved to B2.
$ws.Range("B2").Select()
